$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Rspo3 -> Sdc4 -> ECs
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Rspo3"
$ws.Cells.Item(2,3).Value = "Sdc4"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 2.507621333333333
$ws.Cells.Item(2,8).Value = 7.522864
$ws.Cells.Item(2,9).Value = 0.9939780200440224
$ws.Cells.Item(2,10).Value = 0.9939780200440224
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 2.300909333333333
$ws.Cells.Item(2,14).Value = 6.902728
$ws.Cells.Item(2,15).Value = 0.03776979643482627
$ws.Cells.Item(2,16).Value = 0.03776979643482627
$ws.Cells.Item(2,17).Value = 5.769809330332444
$ws.Cells.Item(2,18).Value = 51.928283972992
$ws.Cells.Item(2,19).Value = 0.0375423474777544
$ws.Cells.Item(2,20).Value = 0.03754234747775439

# Row 3: FAPs -> Rspo3 -> Sdc4 -> FAPs
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Rspo3"
$ws.Cells.Item(3,3).Value = "Sdc4"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 2.507621333333333
$ws.Cells.Item(3,8).Value = 7.522864
$ws.Cells.Item(3,9).Value = 0.9939780200440224
$ws.Cells.Item(3,10).Value = 0.9939780200440224
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 16.96312166666667
$ws.Cells.Item(3,14).Value = 50.889365
$ws.Cells.Item(3,15).Value = 0.2784523679257784
$ws.Cells.Item(3,16).Value = 0.2784523679257784
$ws.Cells.Item(3,17).Value = 42.53708577126222
$ws.Cells.Item(3,18).Value = 382.83377194136
$ws.Cells.Item(3,19).Value = 0.2767755333474349
$ws.Cells.Item(3,20).Value = 0.2767755333474349

# Row 4: FAPs -> Rspo3 -> Sdc4 -> sCs
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Rspo3"
$ws.Cells.Item(4,3).Value = "Sdc4"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 2.507621333333333
$ws.Cells.Item(4,8).Value = 7.522864
$ws.Cells.Item(4,9).Value = 0.9939780200440224
$ws.Cells.Item(4,10).Value = 0.9939780200440224
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 41.655263
$ws.Cells.Item(4,14).Value = 124.965789
$ws.Cells.Item(4,15).Value = 0.6837778356393953
$ws.Cells.Item(4,16).Value = 0.6837778356393953
$ws.Cells.Item(4,17).Value = 104.4556261444107
$ws.Cells.Item(4,18).Value = 940.100635299696
$ws.Cells.Item(4,19).Value = 0.6796601392188331
$ws.Cells.Item(4,20).Value = 0.6796601392188331

# Row 5 (new): sCs -> Rspo3 -> Sdc4 -> ECs
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "Rspo3"
$ws.Cells.Item(5,3).Value = "Sdc4"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.01519233333333333
$ws.Cells.Item(5,8).Value = 0.045577
$ws.Cells.Item(5,9).Value = 0.006021979955977724
$ws.Cells.Item(5,10).Value = 0.006021979955977723
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 2.300909333333333
$ws.Cells.Item(5,14).Value = 6.902728
$ws.Cells.Item(5,15).Value = 0.03776979643482627
$ws.Cells.Item(5,16).Value = 0.03776979643482627
$ws.Cells.Item(5,17).Value = 0.03495618156177778
$ws.Cells.Item(5,18).Value = 0.314605634056
$ws.Cells.Item(5,19).Value = 0.0002274489570718827
$ws.Cells.Item(5,20).Value = 0.0002274489570718827

# Row 6 (new): sCs -> Rspo3 -> Sdc4 -> FAPs
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Rspo3"
$ws.Cells.Item(6,3).Value = "Sdc4"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.01519233333333333
$ws.Cells.Item(6,8).Value = 0.045577
$ws.Cells.Item(6,9).Value = 0.006021979955977724
$ws.Cells.Item(6,10).Value = 0.006021979955977723
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 16.96312166666667
$ws.Cells.Item(6,14).Value = 50.889365
$ws.Cells.Item(6,15).Value = 0.2784523679257784
$ws.Cells.Item(6,16).Value = 0.2784523679257784
$ws.Cells.Item(6,17).Value = 0.2577093987338889
$ws.Cells.Item(6,18).Value = 2.319384588605
$ws.Cells.Item(6,19).Value = 0.001676834578343572
$ws.Cells.Item(6,20).Value = 0.001676834578343572

# Row 7 (new): sCs -> Rspo3 -> Sdc4 -> sCs
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Rspo3"
$ws.Cells.Item(7,3).Value = "Sdc4"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.01519233333333333
$ws.Cells.Item(7,8).Value = 0.045577
$ws.Cells.Item(7,9).Value = 0.006021979955977724
$ws.Cells.Item(7,10).Value = 0.006021979955977723
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 41.655263
$ws.Cells.Item(7,14).Value = 124.965789
$ws.Cells.Item(7,15).Value = 0.6837778356393953
$ws.Cells.Item(7,16).Value = 0.6837778356393953
$ws.Cells.Item(7,17).Value = 0.6328406405836666
$ws.Cells.Item(7,18).Value = 5.695565765253
$ws.Cells.Item(7,19).Value = 0.00411769642056227
$ws.Cells.Item(7,20).Value = 0.004117696420562269
